$d = $word.ActiveDocument

# 1. Merge "40" + "772" into "40781" (Changeset number update)
$d.Content.Find.Execute("772", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "781", 2)

# 2. Merge "Go to " + " " into "Go to  " (collapse the two runs before the hyperlink)
$d.Content.Find.Execute("Go to  ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Go to  ", 2)
